$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new value looks numeric must be forced to Text so
# Excel does not silently convert them (losing exact string formatting,
# e.g. trailing zeros) the same way the source data keeps them as text.

$ws.Range("D2").Value = '60.153.91'
$ws.Range("E2").Value = '  +3.55%  '
$ws.Range("D3").Value = '2.422.50'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '554.22'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.09%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '137.63'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.32%  '
$ws.Range("E7").Value = '  +0.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.579'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("E9").Value = '  +3.26%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '5.78'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.48%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  -1.88%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '24.94'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.60%  '
$ws.Range("D14").Value = '2.851.85'
$ws.Range("E14").Value = '  +3.21%  '
$ws.Range("D15").Value = '60.068.15'
$ws.Range("E15").Value = '  +3.53%  '
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("D17").Value = '2.418.75'
$ws.Range("E17").Value = '  +4.56%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '11.36'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +6.28%  '
$ws.Range("E19").Value = '  +2.16%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '332.23'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +3.91%  '
$ws.Range("E24").Value = '  +3.56%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '8.60'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.28%  '
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '0.0₃0784'
$ws.Range("E28").Value = '  +6.39%  '
$ws.Range("E29").Value = '  +1.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '169.99'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("E31").Value = '  +1.73%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '18.66'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("E35").Value = '  +5.29%  '
$ws.Range("E36").Value = '  +0.17%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.21'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.22%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.61'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.28%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '39.54'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("E40").Value = '  +10.59%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '313.47'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +8.33%  '
$ws.Range("E42").Value = '  +1.42%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '139.20'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.59%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0962'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.64%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0520'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.92%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.50%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.411'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +7.85%  '
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("E49").Value = '  +1.38%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '17.73'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.90%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '11.05'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.24%  '
